# Lead time + small bug
# - implement lead time (works)
# - bug in cheapest product calculation
# - To do: DEBUG VSS calculation

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "three_scenarios": only the view/selection moved (no data changed).
# ---------------------------------------------------------------------------
$wsThree = $wb.Worksheets.Item("three_scenarios")
$wsThree.Activate()
$wsThree.Range("K3").Select()

# ---------------------------------------------------------------------------
# Sheet "three_scenarios_new": fix the 1/3 -> 1/4 probability bug (now that a
# 4th scenario row is added) and add the new "MMM" (always-include) row used
# for the VSS calculation, plus a note explaining why.
# ---------------------------------------------------------------------------
$wsNew = $wb.Worksheets.Item("three_scenarios_new")
$wsNew.Activate()

$wsNew.Range("C2").Formula = "=1/4"
$wsNew.Range("C3").Formula = "=1/4"
$wsNew.Range("C4").Formula = "=1/4"

$wsNew.Range("A5").Value = 3
$wsNew.Range("B5").Value = "MMM"
$wsNew.Range("C5").Formula = "=1/4"
$wsNew.Range("D5").Value = 1
$wsNew.Range("E5").Value = 1
$wsNew.Range("F5").Value = 1
$wsNew.Range("G5").Value = 1
$wsNew.Range("H5").Value = "base"
$wsNew.Range("I5").Value = "base"
$wsNew.Range("J5").Value = "base"
$wsNew.Range("K5").Value = "base"
$wsNew.Range("M5").Value = 'always include "MMM" -> needed for calculation of VSS'

$wsNew.Range("M5").Select()
